# Fix negative Average Generation Cost values: the "Average Generation
# Cost" column (C) had ended up with the wrong figures (small ratio-like
# numbers, some negative/odd) — they should instead mirror the raw kWh
# figures already present in column B for the same row (rows 2-131).
#
# We force each target cell to Text format ("@") before writing so that
# values containing thousands separators, decimals, or parentheses are
# copied over verbatim (as literal text) instead of being silently
# re-parsed/reformatted as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 131) { $lastRow = 131 }

for ($r = 2; $r -le $lastRow; $r++) {
    $srcCell = $ws.Cells.Item($r, 2)
    $dstCell = $ws.Cells.Item($r, 3)
    $dstCell.NumberFormat = "@"
    $dstCell.Value = $srcCell.Value()
}
